# Updates the "cryptos" price/volume table with refreshed values.
# Note: several "Price" (column D) entries are numeric-looking text
# (e.g. "309.40"); a plain .Value assignment would let Excel coerce
# them to real numbers. To keep them as text (matching the source
# data, which stores thousands-separated / precise price strings),
# we briefly force NumberFormat to Text ("@") before assigning, then
# restore the cell's original (default) style by copying it from an
# adjacent untouched cell in the same row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.927.86'
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").Value = '1.846.41'
$ws.Range("E3").Value = '  +1.24%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.40'
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4768'
$ws.Range("D7").Style = $ws.Range("C7").Style
$ws.Range("E7").Value = '  +2.85%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3666'
$ws.Range("D8").Style = $ws.Range("C8").Style
$ws.Range("E8").Value = '  +1.73%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07205'
$ws.Range("D9").Style = $ws.Range("C9").Style
$ws.Range("E9").Value = '  +1.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9271'
$ws.Range("D10").Style = $ws.Range("C10").Style
$ws.Range("E10").Value = '  +2.98%  '
$ws.Range("E11").Value = '  +1.86%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.939.32'
$ws.Range("E12").Value = '  +7.17%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07687'
$ws.Range("D13").Style = $ws.Range("C13").Style
$ws.Range("E13").Value = '  -1.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.314'
$ws.Range("D14").Style = $ws.Range("C14").Style
$ws.Range("E14").Value = '  +1.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.402'
$ws.Range("D15").Style = $ws.Range("C15").Style
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.81'
$ws.Range("D16").Style = $ws.Range("C16").Style
$ws.Range("E16").Value = '  +1.44%  '
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008626'
$ws.Range("D18").Style = $ws.Range("C18").Style
$ws.Range("E18").Value = '  +0.69%  '
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").Value = '26.952.27'
$ws.Range("E20").Value = '  +1.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.049'
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = '  +0.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.63'
$ws.Range("D23").Style = $ws.Range("C23").Style
$ws.Range("E23").Value = '  +0.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.919'
$ws.Range("D24").Style = $ws.Range("C24").Style
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.44'
$ws.Range("D25").Style = $ws.Range("C25").Style
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.15'
$ws.Range("D26").Style = $ws.Range("C26").Style
$ws.Range("E26").Value = '  +1.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.995'
$ws.Range("D27").Style = $ws.Range("C27").Style
$ws.Range("E27").Value = '  +1.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '114.24'
$ws.Range("D28").Style = $ws.Range("C28").Style
$ws.Range("E28").Value = '  +0.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.919'
$ws.Range("D29").Style = $ws.Range("C29").Style
$ws.Range("E29").Value = '  +2.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08879'
$ws.Range("D30").Style = $ws.Range("C30").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.312'
$ws.Range("D31").Style = $ws.Range("C31").Style
$ws.Range("E31").Value = '  +5.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.175'
$ws.Range("D32").Style = $ws.Range("C32").Style
$ws.Range("E32").Value = '  +3.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7439'
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").Value = '  +1.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.484'
$ws.Range("D34").Style = $ws.Range("C34").Style
$ws.Range("E34").Value = '  +1.11%  '
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.105'
$ws.Range("D36").Style = $ws.Range("C36").Style
$ws.Range("E36").Value = '  +2.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01954'
$ws.Range("D37").Style = $ws.Range("C37").Style
$ws.Range("E37").Value = '  +1.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05260'
$ws.Range("D38").Style = $ws.Range("C38").Style
$ws.Range("E38").Value = '  +2.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.972'
$ws.Range("D39").Style = $ws.Range("C39").Style
$ws.Range("E39").Value = '  +1.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5195'
$ws.Range("D40").Style = $ws.Range("C40").Style
$ws.Range("E40").Value = '  +2.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.946'
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = '  +0.91%  '
$ws.Range("E42").Value = '  +0.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.208'
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("E43").Value = '  +2.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.56'
$ws.Range("D44").Style = $ws.Range("C44").Style
$ws.Range("E44").Value = '  +5.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4715'
$ws.Range("D45").Style = $ws.Range("C45").Style
$ws.Range("E45").Value = '  +1.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.009'
$ws.Range("D46").Style = $ws.Range("C46").Style
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.15'
$ws.Range("D47").Style = $ws.Range("C47").Style
$ws.Range("E47").Value = '  +3.08%  '
$ws.Range("E48").Value = '  +2.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '65.69'
$ws.Range("D49").Style = $ws.Range("C49").Style
$ws.Range("E49").Value = '  +3.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06021'
$ws.Range("D50").Style = $ws.Range("C50").Style
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8847'
$ws.Range("D51").Style = $ws.Range("C51").Style
$ws.Range("E51").Value = '  +3.88%  '
